{"js": "// Update the 25 \"NNN\u00d7N=\" practice-problem cells in the first (and only)\n// table. Cells are addressed positionally (row, column) rather than by\n// searching for the old text, because one of the new values (\"547\u00d77=\",\n// the replacement for row 20 col 2's old \"904\u00d77=\") is identical to an\n// old value used earlier in the document (row 15 col 4's old \"547\u00d77=\",\n// which becomes \"717\u00d72=\"). Addressing by cell coordinates keeps every\n// write independent of the others, so there is no risk of a later\n// replacement accidentally matching text that an earlier replacement\n// just inserted.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based row index -> new values for that row's 5 cells (left to right).\nconst rowUpdates = {\n  0: [\"149\u00d75=\", \"202\u00d75=\", \"430\u00d76=\", \"919\u00d78=\", \"339\u00d78=\"],\n  4: [\"328\u00d77=\", \"527\u00d78=\", \"917\u00d79=\", \"208\u00d72=\", \"204\u00d79=\"],\n  9: [\"260\u00d72=\", \"929\u00d77=\", \"212\u00d72=\", \"539\u00d79=\", \"396\u00d76=\"],\n  14: [\"666\u00d79=\", \"637\u00d72=\", \"650\u00d77=\", \"717\u00d72=\", \"246\u00d76=\"],\n  19: [\"103\u00d76=\", \"547\u00d77=\", \"712\u00d77=\", \"272\u00d74=\", \"902\u00d75=\"],\n};\n\nfor (const rowKey of Object.keys(rowUpdates)) {\n  const row = parseInt(rowKey, 10);\n  const values = rowUpdates[rowKey];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"NNN\u00d7N=\" practice-problem cells in the first (and only)\n# table. Cells are addressed positionally (row, column) rather than by\n# searching for the old text, because one of the new values (\"547\u00d77=\",\n# the replacement for row 20 col 2's old \"904\u00d77=\") is identical to an\n# old value used earlier in the document (row 15 col 4's old \"547\u00d77=\",\n# which becomes \"717\u00d72=\"). Addressing by cell coordinates keeps every\n# write independent of the others, so there is no risk of a later\n# replacement accidentally matching text that an earlier replacement\n# just inserted.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# row -> list of 5 new values (left to right) for that row\n$rowUpdates = @{\n    1  = @(\"149\u00d75=\", \"202\u00d75=\", \"430\u00d76=\", \"919\u00d78=\", \"339\u00d78=\")\n    5  = @(\"328\u00d77=\", \"527\u00d78=\", \"917\u00d79=\", \"208\u00d72=\", \"204\u00d79=\")\n    10 = @(\"260\u00d72=\", \"929\u00d77=\", \"212\u00d72=\", \"539\u00d79=\", \"396\u00d76=\")\n    15 = @(\"666\u00d79=\", \"637\u00d72=\", \"650\u00d77=\", \"717\u00d72=\", \"246\u00d76=\")\n    20 = @(\"103\u00d76=\", \"547\u00d77=\", \"712\u00d77=\", \"272\u00d74=\", \"902\u00d75=\")\n}\n\nforeach ($row in $rowUpdates.Keys) {\n    $values = $rowUpdates[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $cell = $tbl.Cell($row, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
